$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize the query table / ListObject to make room for the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B366"))

# --- Append the 8 new municipality rows pulled in by the refreshed query ---
$ws.Range("A359").Value = 5007
$ws.Range("B359").Value = "Namsos"

$ws.Range("A360").Value = 1826
$ws.Range("B360").Value = "Hattfjelldal"

$ws.Range("A361").Value = 1841
$ws.Range("B361").Value = "Fauske Fuosko"

$ws.Range("A362").Value = 1853
$ws.Range("B362").Value = "Evenes"

$ws.Range("A363").Value = 1875
$ws.Range("B363").Value = "Hamarøy Hábmer"

$ws.Range("A364").Value = 5406
$ws.Range("B364").Value = "Hammerfest"

$ws.Range("A365").Value = 5412
$ws.Range("B365").Value = "Tjeldsund"

$ws.Range("A366").Value = 5436
$ws.Range("B366").Value = "Porsanger Porsáŋgu Porsanki"

# New codes come in as plain numbers (the refreshed source typed them that
# way), so give column A an explicit General format; column B picks up the
# same explicit formatting from row 360 onward (row 359 was the first one
# typed in and kept the table's inherited default format).
$ws.Range("A359:A366").NumberFormat = "General"
$ws.Range("B360:B366").NumberFormat = "General"

# --- Keep the hidden query-table defined name in sync with the new extent ---
$wb.Names.Item("EksterneData_2").RefersTo = "=Ark1!`$A`$1:`$B`$366"

# --- Move the selection to the first blank row below the refreshed data ---
$ws.Range("A367").Select() | Out-Null
